$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 601
$ws.Range("F5").Value = 1422
$ws.Range("F6").Value = 731
$ws.Range("F9").Value = 444
$ws.Range("F10").Value = 6643
$ws.Range("F11").Value = 131
$ws.Range("F14").Value = 4895
$ws.Range("F16").Value = 6061
$ws.Range("F17").Value = 7804
$ws.Range("F18").Value = 152
$ws.Range("F21").Value = 4128
$ws.Range("F22").Value = 600
$ws.Range("F23").Value = 64
$ws.Range("F27").Value = 1094
$ws.Range("F29").Value = 1537
$ws.Range("F30").Value = 602
$ws.Range("F32").Value = 1738
$ws.Range("F34").Value = 2010
$ws.Range("F35").Value = 248
$ws.Range("F36").Value = 63
$ws.Range("F37").Value = 1299
$ws.Range("F39").Value = 722
$ws.Range("F41").Value = 3801
$ws.Range("F45").Value = 467
$ws.Range("F47").Value = 32
$ws.Range("F48").Value = 115

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4670

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4670
$ws.Range("F9").Value = 601
$ws.Range("F11").Value = 1422
$ws.Range("F12").Value = 731
$ws.Range("F14").Value = 444
$ws.Range("F15").Value = 6643
$ws.Range("F18").Value = 4895
$ws.Range("F19").Value = 6061
$ws.Range("F20").Value = 6061
$ws.Range("F21").Value = 7804
$ws.Range("F22").Value = 152
$ws.Range("F25").Value = 4128
$ws.Range("F26").Value = 600
$ws.Range("F27").Value = 64
$ws.Range("F30").Value = 1094
$ws.Range("F31").Value = 1537
$ws.Range("F32").Value = 602
$ws.Range("F34").Value = 1738
$ws.Range("F36").Value = 2010
$ws.Range("F41").Value = 722
$ws.Range("F45").Value = 3801
$ws.Range("F47").Value = 371
$ws.Range("F48").Value = 32
$ws.Range("F49").Value = 115
